$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6+ down by one.
$ws.Rows.Item(6).Insert()

# Fill in the newly inserted row 6 with the new agenda entry.
$ws.Cells.Item(6, 1).Value = "Giovani"
$ws.Cells.Item(6, 2).Value = "0446"
$ws.Cells.Item(6, 3).Value = "Insecta"
$ws.Cells.Item(6, 4).Value = "Linha telefônica, aparentemente sem comunicação de alarmes."
$ws.Cells.Item(6, 7).Value = "Pendente"

# The newly inserted row picks up default formatting; copy the (now
# shifted-down) formatting that used to belong to row 6 back onto it,
# matching Excel's normal "insert row" behaviour of carrying the format
# of the row above down into the new row. Done after setting the values
# so the quote-prefix / number-format styling isn't clobbered again.
$ws.Range("A7:I7").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to match the recorded cursor position.
$ws.Range("H6").Select()
